$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.524.32"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "2.240.69"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").Value = "'306.36"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "'94.88"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").Value = "'34.92"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "2.278.87"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").Value = "'0.833"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "'13.55"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "44.256.08"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "0.0₃0952"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").Value = "'11.91"
$ws.Range("E20").Value = "  -1.69%  "

$ws.Range("D21").Value = "'65.51"
$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").Value = "'237.38"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").Value = "'2.95"
$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").Value = "'37.76"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").Value = "'9.77"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").Value = "'5.97"
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("D30").Value = "'19.92"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").Value = "'152.62"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").Value = "'0.0796"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +2.26%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.02"
$ws.Range("E35").Value = "  -6.84%  "

$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("E37").Value = "  +2.05%  "

$ws.Range("D38").Value = "'15.00"
$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("D40").Value = "'3.77"
$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("D41").Value = "'0.0300"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").Value = "1.798.16"
$ws.Range("E43").Value = "  +4.27%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.69"
$ws.Range("E44").Value = "  +12.49%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.191"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").Value = "'78.61"
$ws.Range("E46").Value = "  -7.97%  "

$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'4.91"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'98.71"
$ws.Range("E48").Value = "  -1.43%  "

$ws.Range("D49").Value = "'69.92"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("D50").Value = "'8.09"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").Value = "'54.43"
$ws.Range("E51").Value = "  +0.30%  "
